$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.018.66"
$ws.Range("E2").Value = "  -0.96%  "

$ws.Range("D3").Value = "2.947.33"
$ws.Range("E3").Value = "  +0.56%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "378.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.99%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.37"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.537"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.16%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.586"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.14%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.28"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.95%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0836"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.24%  "

$ws.Range("D13").Value = "3.410.72"
$ws.Range("E13").Value = "  +0.64%  "

$ws.Range("E14").Value = "  -2.52%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.37"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.00%  "

$ws.Range("D16").Value = "2.941.27"
$ws.Range("E16").Value = "  +0.60%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.982"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.64%  "

$ws.Range("D18").Value = "51.071.73"
$ws.Range("E18").Value = "  -0.81%  "

$ws.Range("E19").Value = "  -6.81%  "

$ws.Range("E20").Value = "  -0.87%  "

$ws.Range("E21").Value = "  -3.86%  "

$ws.Range("E22").Value = "  +0.15%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.35"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.03%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "260.61"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.65%  "

$ws.Range("E25").Value = "  +3.28%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.16"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +11.51%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.56"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +8.90%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.10"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.86%  "

$ws.Range("E29").Value = "  -0.02%  "

$ws.Range("E30").Value = "  +10.55%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.165"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.60%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "25.65"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.39%  "

$ws.Range("E33").Value = "  -0.59%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "50.45"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.20%  "

$ws.Range("E35").Value = "  -2.99%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "33.45"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.73%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0442"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.32%  "

$ws.Range("E38").Value = "  -0.03%  "

$ws.Range("E39").Value = "  -1.75%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.86"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.75%  "

$ws.Range("E41").Value = "  +0.26%  "

$ws.Range("E42").Value = "  -2.56%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.77"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.88%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "121.41"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.76%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.01"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.65%  "

$ws.Range("E46").Value = "  +0.34%  "

$ws.Range("E47").Value = "  -0.33%  "

$ws.Range("D49").Value = "1.997.02"
$ws.Range("E49").Value = "  -1.27%  "

$ws.Range("E50").Value = "  +0.93%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0330"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.85%  "
